# Generate Report for handoff
#
# The zh-cn / de-de handoff attempt failed during transform, so the
# per-language status rows need to reflect that: the status flips from
# "Ready for handoff" to "Handoff transform failed", the handoff file
# link/date are cleared back to defaults, and the handoff reason flips
# from "Include" to "Ignored".

$wb = $excel.ActiveWorkbook

# The Overview sheet mirrors each language's status in its own column, so
# it needs the same status text update.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"

$sheetNames = @("zh-cn", "de-de")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column (B2): "Ready for handoff" -> "Handoff transform failed"
    $ws.Range("B2").Value = "Handoff transform failed"

    # Latest Handoff File (C2): remove the hyperlink + its display text,
    # the handoff artifact no longer exists since the transform failed.
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Address -like "*.xlf") {
            $hl.Delete()
        }
    }
    $ws.Range("C2").ClearContents()
    $ws.Range("C2").ClearFormats()

    # Latest Handoff Datetime (D2): reset to the default "unset" datetime.
    $ws.Range("D2").Value = "0001-01-01 00:00:00"

    # Handoff Reason (H2): "Include" -> "Ignored"
    $ws.Range("H2").Value = "Ignored"
}
